# Update countries & provincias Spain
# Refresh COVID data table ("Pais" sheet): updated timestamp, refreshed
# case numbers for a handful of countries, and a couple of countries
# (Marruecos, Sri Lanka) that moved a few rows up in the ranking now that
# their totals overtook their neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 12:52"

# Suiza (row 18) - refreshed numbers, same country/position
$ws.Range("B18").Value = 27944
$ws.Range("C18").Value = 204
$ws.Range("D18").Value = 17800
$ws.Range("E18").Value = 8738
$ws.Range("F18").Value = 386
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 1406

# Marruecos moves up from row 57 to row 55 (ahead of Banglades/Argentina),
# with refreshed totals; Banglades and Argentina each shift down one row,
# keeping their previous totals.
$ws.Range("A55").Value = "Marruecos"
$ws.Range("B55").Value = 2990
$ws.Range("C55").Value = 135
$ws.Range("D55").Value = 340
$ws.Range("E55").Value = 2507
$ws.Range("F55").Value = 1
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 143

$ws.Range("A56").Value = "Banglades"
$ws.Range("B56").Value = 2948
$ws.Range("C56").Value = 492
$ws.Range("D56").Value = 85
$ws.Range("E56").Value = 2762
$ws.Range("F56").Value = 1
$ws.Range("G56").Value = 10
$ws.Range("H56").Value = 101

$ws.Range("A57").Value = "Argentina"
$ws.Range("B57").Value = 2941
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 709
$ws.Range("E57").Value = 2098
$ws.Range("F57").Value = 123
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 134

# Moldavia (row 60) - refreshed numbers, same country/position
$ws.Range("E60").Value = 1947
$ws.Range("F60").Value = 212
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 68

# Kazajistan (row 67) - refreshed numbers, same country/position
$ws.Range("D67").Value = 428
$ws.Range("E67").Value = 1310

# Uzbekistan (row 68) - refreshed numbers, same country/position
$ws.Range("D68").Value = 240
$ws.Range("E68").Value = 1337

# Bosnia y Herzegovina (row 77) - refreshed numbers, same country/position
$ws.Range("B77").Value = 1309
$ws.Range("C77").Value = 24
$ws.Range("D77").Value = 381
$ws.Range("E77").Value = 879
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 49

# Libano (row 93) - refreshed numbers, same country/position
$ws.Range("B93").Value = 677
$ws.Range("C93").Value = 4
$ws.Range("E93").Value = 554

# Malta (row 106) - refreshed numbers, same country/position
$ws.Range("B106").Value = 431
$ws.Range("C106").Value = 4
$ws.Range("D106").Value = 126
$ws.Range("E106").Value = 302

# Senegal (row 111) - refreshed numbers, same country/position
$ws.Range("B111").Value = 377
$ws.Range("C111").Value = 10
$ws.Range("D111").Value = 235
$ws.Range("E111").Value = 137
$ws.Range("G111").Value = 2
$ws.Range("H111").Value = 5

# Sri Lanka moves up from row 116 to row 115 (ahead of Isla de Man), with
# refreshed totals; Isla de Man shifts down one row, keeping its previous
# totals.
$ws.Range("A115").Value = "Sri Lanka"
$ws.Range("B115").Value = 303
$ws.Range("C115").Value = 32
$ws.Range("D115").Value = 97
$ws.Range("E115").Value = 199
$ws.Range("F115").Value = 1
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 7

$ws.Range("A116").Value = "Isla de Man"
$ws.Range("B116").Value = 298
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 193
$ws.Range("E116").Value = 99
$ws.Range("F116").Value = 10
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 6
